# ---------------------------------------------------------------------------
# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# 1. "ODI Batting" sheet:
#    - MATCH_CARD_LINK column (D) renamed to MATCH_CODE, and every URL
#      value replaced by just the trailing MatchCode number (as text).
#    - A handful of rows had a stray empty INNING_NUMBER (B) cell; those
#      are cleared out entirely.
# 2. A new "Player Info" sheet is inserted in front of "ODI Batting" with
#    the player's bio info.
# 3. A new "ODI Batting Extra" sheet is appended after "ODI Batting" with
#    extra per-match batting stats for the most recent matches.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$odi = $wb.Worksheets.Item("ODI Batting")

# --- 1a. Rename MATCH_CARD_LINK -> MATCH_CODE, strip the URL down to the
#         bare MatchCode number (kept as text, matching the other text
#         columns in the sheet). ------------------------------------------
$odi.Range("D1").Value = "MATCH_CODE"

$lastRow = $odi.Cells.Item($odi.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $odi.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url -match "MatchCode=(\d+)") {
        $cell.Value = "'" + $matches[1]
    }
}

# --- 1b. Clear the stray empty INNING_NUMBER (column B) cells. -----------
$emptyBRows = @(18, 30, 37, 39, 54, 55, 56, 63, 67, 68, 75, 86, 89, 90)
foreach ($r in $emptyBRows) {
    $odi.Cells.Item($r, 2).ClearContents()
}

# --- 2. New "Player Info" sheet, inserted before "ODI Batting". ----------
$playerInfo = $wb.Worksheets.Add($odi)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$playerInfo.Range("A2").Value = "'3910"
$playerInfo.Range("B2").Value = "Matthew Scott Wade"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Does Not Bowl | Unknown"

$odi.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. New "ODI Batting Extra" sheet, appended after "ODI Batting". -----
# (Worksheets.Add ignores the "After" slot in this host, so add it anywhere
# and then explicitly Move it into place.)
$extra = $wb.Worksheets.Add($odi)
$extra.Name = "ODI Batting Extra"
$odiRef = $wb.Worksheets.Item("ODI Batting")
$extra.Move($null, $odiRef)

$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

$odi.Range("A1").Copy()
$extra.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("'3947", $null, $null, $null, $null, "NO"),
    @("'3950", $null, $null, $null, $null, "NO"),
    @("'3966", $null, $null, $null, $null, "NO"),
    @("'3967", 6, "'1", "'1", "'2.91%", "NO"),
    @("'3968", $null, $null, $null, $null, "NO"),
    @("'3972", 7, "'7", "'2", "'37.31%", "YES"),
    @("'3973", $null, $null, $null, $null, "NO"),
    @("'3975", 7, $null, $null, $null, "NO"),
    @("'3977", 6, "'0", "'0", "'1.42%", "NO"),
    @("'3981", 5, "'1", "'0", "'2.17%", "NO"),
    @("'4032", 7, $null, $null, $null, "NO"),
    @("'4035", 7, $null, $null, $null, "NO"),
    @("'4041", $null, $null, $null, $null, "NO"),
    @("'4067", 7, "'1", "'0", "'6.57%", "NO"),
    @("'4069", 7, "'0", "'0", "'0.99%", "NO"),
    @("'4074", $null, $null, $null, $null, "NO"),
    @("'4076", 7, "'0", "'1", "'8.26%", "NO"),
    @("'4483", 7, "'0", "'0", "'1.19%", "NO"),
    @("'4484", 7, "'2", "'0", "'19.25%", "NO"),
    @("'4486", $null, $null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraRows) {
    $extra.Cells.Item($r, 1).Value = $row[0]
    if ($null -ne $row[1]) {
        $extra.Cells.Item($r, 2).Value = $row[1]
    }
    if ($null -ne $row[2]) {
        $extra.Cells.Item($r, 3).Value = $row[2]
    }
    if ($null -ne $row[3]) {
        $extra.Cells.Item($r, 4).Value = $row[3]
    }
    if ($null -ne $row[4]) {
        $extra.Cells.Item($r, 5).Value = $row[4]
    }
    $extra.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$playerInfo.Range("A1").Select()
